$d = $word.ActiveDocument

# ------------------------------------------------------------------
# 1) Locate the paragraph that ends with "...«unassigned»)." (the
#    "teacher" bullet point), change its final "." into ";" and
#    insert a brand-new bullet (same list/style) right after it
#    describing the administrator's per-class table view.
# ------------------------------------------------------------------
$teacherIdx = 0
$i = 0
foreach ($p in $d.Paragraphs) {
    $i = $i + 1
    if ($p.Range.Text -like "*unassigned*") {
        $teacherIdx = $i
    }
}

$teacherPara = $d.Paragraphs.Item($teacherIdx)
$teacherRange = $teacherPara.Range
# Last real character of the paragraph sits right before the paragraph
# mark, i.e. between (End-2) and (End-1).
$period = $d.Range($teacherRange.End - 2, $teacherRange.End - 1)
$period.Text = ";"

# Re-fetch the paragraph range (offsets/handles may shift after edits)
$teacherPara = $d.Paragraphs.Item($teacherIdx)
$teacherPara.Range.InsertParagraphAfter()

$adminClassPara = $d.Paragraphs.Item($teacherIdx + 1)
$adminClassPara.Range.InsertBefore("Для администратора реализовано отображение информации по всем классам школы в виде таблицы, где напротив имени учащегося расположена кнопка, нажав на которую, администратор может отдать ученика родителям, если статус его заявки «выпущен», и создать заявку, если вместо статуса стоит прочерк.")

# ------------------------------------------------------------------
# 2) Append a new bullet (same list/style as the other client-side
#    bullets) after the very last paragraph of the document, which
#    describes the administrator grouped-info interface.
# ------------------------------------------------------------------
$lastPara = $d.Paragraphs.Item($d.Paragraphs.Count)
$lastPara.Range.InsertParagraphAfter()

$adminGroupPara = $d.Paragraphs.Item($d.Paragraphs.Count)
$adminGroupPara.Range.InsertBefore("Отображение сгруппированной информации о времени выхода и статусе заявки всех учеников школы в интерфейсе администратора с возможностью создать новую заявку (см. пункт 5).")
